# chore: publish terminology IG 2.0.2
# Update the metadata table on the "Metadata" sheet:
#   Version:      1.8.1 -> 1.8.2
#   Status:       draft -> active
#   Experimental: true  -> (blank)
#   Date:         2024-01-18 -> 2025-11-18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version
$ws.Range("B3").Value = "1.8.2"

# Status
$ws.Range("B6").Value = "active"

# Experimental (cleared, as the resource is no longer marked experimental)
$ws.Range("B7").ClearContents()

# Date - write as plain text ("2025-11-18") rather than letting Excel
# auto-convert the string into a date serial number. We stage the text
# in an unused helper cell that is explicitly formatted as Text, copy
# it, and paste-special (values only) into the target cell so the
# target keeps its original style/number format untouched.
$helper = $ws.Range("D1")
$helper.NumberFormat = "@"
$helper.Value = "2025-11-18"
$helper.Copy() | Out-Null
$ws.Range("B8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
$helper.Clear() | Out-Null
$excel.CutCopyMode = $false
